# "ore di lavoro aggiornate" - update Ore_di_lavoro.xlsx / Foglio1
# Add the hours logged for 28/09/2022 (row 86, serial date 44690) plus
# the two new note lines describing what was worked on that day.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Hours worked that day (column G = Jacopo)
$ws.Range("G86").Value = 4

# Notes entered under that day's row, in the order they were typed so the
# shared-string table picks up "coding , testing " before
# "taiga, documentazione".
$ws.Range("G88").Value = "coding , testing "
$ws.Range("G87").Value = "taiga, documentazione"
